# Update to the existing API example file to match the deployment changes.
# Adds new error-row documentation to three sheets:
#   - PMOD(Open and Distribute)      rows 6-12  (dimension A1:G9 -> A1:G12)
#   - Domestic Labels Outbound       rows 206-209 (dimension A1:G205 -> A1:G209)
#   - Domestic Labels Returns        rows 182-186 (dimension A1:G181 -> A1:G186)

$wb = $excel.ActiveWorkbook

$pmod       = $wb.Worksheets.Item("PMOD(Open and Distribute)")
$outbound   = $wb.Worksheets.Item("Domestic Labels Outbound")
$returns    = $wb.Worksheets.Item("Domestic Labels Returns ")
$containers = $wb.Worksheets.Item("Containers")
$genErrors  = $wb.Worksheets.Item("General Errors")

# ---------------------------------------------------------------------------
# 1) PMOD(Open and Distribute) sheet
# ---------------------------------------------------------------------------

# Row 5, column D loses its (redundant) fill-flagged wrap style, switching to
# the plain wrap-text style used elsewhere on the sheet.
$containers.Range("G2").Copy()
$pmod.Range("D5").PasteSpecial(-4122)

# Build the formatting for the new rows 6-12 from existing templates, then
# overwrite the 3 previously-blank placeholder rows (7-9) and extend with
# fresh rows (10-12).
$pmod.Range("A4:G4").Copy()
$pmod.Range("A6:G12").PasteSpecial(-4122)

$pmod.Range("D3").Copy()
$pmod.Range("D8:D12").PasteSpecial(-4122)

$containers.Range("G2").Copy()
$pmod.Range("G6:G12").PasteSpecial(-4122)

$pmod.Range("A6:A12").Value = "PMOD"
$pmod.Range("B6:B12").Value = 400
$pmod.Range("C6:C12").Value = "BAD_REQUEST"

$pmod.Range("D6").Value = "https://api.usps.com/pmod/v3/pmod/manifest/{container-id}"
$pmod.Range("D7").Value = "https://api.usps.com/pmod/v3/pmod/manifest/{container-id}"
$pmod.Range("D8").Value = "https://api.usps.com/pmod/v3/pmod/create"
$pmod.Range("D9").Value = "https://api.usps.com/pmod/v3/pmod/create"
$pmod.Range("D10").Value = "https://api.usps.com/pmod/v3/pmod/create"
$pmod.Range("D11").Value = "https://api.usps.com/pmod/v3/pmod/create"
$pmod.Range("D12").Value = "https://api.usps.com/pmod/v3/pmod/create"

$pmod.Range("E6").Value = 160344
$pmod.Range("E7").Value = 160345
$pmod.Range("E8").Value = 160350
$pmod.Range("E9").Value = 160352
$pmod.Range("E10").Value = 160353
$pmod.Range("E11").Value = 160354
$pmod.Range("E12").Value = 160355

$pmod.Range("F6").Value = "containerId"
$pmod.Range("F7").Value = "containerId"
$pmod.Range("F8").Value = "trackingNumbers"
$pmod.Range("F9").Value = "trackingNumbers"
$pmod.Range("F10").Value = "trackingNumbers"
$pmod.Range("F11").Value = "trackingNumbers"
$pmod.Range("F12").Value = "trackingNumbers"

$pmod.Range("G6").Value = "Container does not have any packages associated with it"
$pmod.Range("G7").Value = "Container is not a PMOD container"
$pmod.Range("G8").Value = "Unable to find valid label with trackingNumber: <trackingNumber>"
$pmod.Range("G9").Value = "Unable to add Priority Mail Express label <barcode> to Priority Mail Express Open Distribute container"
$pmod.Range("G10").Value = "Unable to add hazmat label <barcode> to Priority Mail Express Open Distribute container"
$pmod.Range("G11").Value = "Label <barcode> with status <status> is unable to be added to container"
$pmod.Range("G12").Value = "Unable to add label <barcode> with mail class <mailClass> to Open Distribute Container"

# View state: scrolled one column right, with a new active selection.
$pmod.Application.ActiveWindow.ScrollColumn = 2
[void]$pmod.Range("G20").Select()

# ---------------------------------------------------------------------------
# 2) Domestic Labels Outbound sheet
# ---------------------------------------------------------------------------

# Rows 206-209 share the same formatting as row 50 (and many others).
$outbound.Range("A50:G50").Copy()
$outbound.Range("A206:G209").PasteSpecial(-4122)

$outbound.Range("A206:A209").Value = "Domestic Labels Outbound"
$outbound.Range("B206:B209").Value = 400
$outbound.Range("C206:C209").Value = "BAD_REQUEST"
$outbound.Range("D206:D209").Value = "https://api.usps.com/labels/v3/label"
$outbound.Range("F206:F209").Value = "imageInfo.brandingImageFormat"

$outbound.Range("E206").Value = 160347
$outbound.Range("E207").Value = 160348
$outbound.Range("E208").Value = 160349
$outbound.Range("E209").Value = 160351

$outbound.Range("G206").Value = "TWO_SQUARES' is supported on the following 'imageInfo.labelType' values: ['4X6LABEL', '4X5LABEL']"
$outbound.Range("G207").Value = "ONE_SQUARE' is supported on the following 'imageInfo.labelType' values: ['4X6LABEL', '4X5LABEL', and '2X7LABEL']"
$outbound.Range("G208").Value = "RECTANGLE' is supported on the following 'imageInfo.labelType' values: ['4X6LABEL', '4X5LABEL']"
$outbound.Range("G209").Value = "RECTANGLE' is only supported for Return Labels with 'imageInfo.labelType' of '4X5LABEL'"

[void]$outbound.Range("D198").Select()

# ---------------------------------------------------------------------------
# 3) Domestic Labels Returns sheet
# ---------------------------------------------------------------------------

# Rows 182-186 reuse row 181's A/B/C formatting, with D/E/F/G patched to the
# styles actually used for this block.
$returns.Range("A181:G181").Copy()
$returns.Range("A182:G186").PasteSpecial(-4122)
$returns.Range("G182:G186").ClearFormats()

$containers.Range("A2").Copy()
$returns.Range("D182:D186").PasteSpecial(-4122)

$genErrors.Range("A2").Copy()
$returns.Range("E182:E186").PasteSpecial(-4122)

$containers.Range("G2").Copy()
$returns.Range("F182:F186").PasteSpecial(-4122)

$returns.Range("A182:A186").Value = "Domestic Labels Outbound"
$returns.Range("B182:B186").Value = 400
$returns.Range("C182:C186").Value = "BAD_REQUEST"
$returns.Range("D182:D186").Value = "https://api.usps.com/labels/v3/return-label"

$returns.Range("E182").Value = 160346
$returns.Range("E183").Value = 160347
$returns.Range("E184").Value = 160348
$returns.Range("E185").Value = 160349
$returns.Range("E186").Value = 160351

$returns.Range("F182").Value = "imageInfo.labelType"
$returns.Range("F183").Value = "imageInfo.brandingImageFormat"
$returns.Range("F184").Value = "imageInfo.brandingImageFormat"
$returns.Range("F185").Value = "imageInfo.brandingImageFormat"
$returns.Range("F186").Value = "imageInfo.brandingImageFormat"

$returns.Range("G182").Value = "4X4LABEL is not supported for this request"
$returns.Range("G183").Value = "TWO_SQUARES' is supported on the following 'imageInfo.labelType' values: ['4X6LABEL', '4X5LABEL']"
$returns.Range("G184").Value = "ONE_SQUARE' is supported on the following 'imageInfo.labelType' values: ['4X6LABEL', '4X5LABEL', and '2X7LABEL']"
$returns.Range("G185").Value = "RECTANGLE' is supported on the following 'imageInfo.labelType' values: ['4X6LABEL', '4X5LABEL']"
$returns.Range("G186").Value = "RECTANGLE' is only supported for Return Labels with 'imageInfo.labelType' of '4X5LABEL'"

[void]$returns.Range("E183:G186").Select()

Write-Output "edit complete"
